$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.042.89"
$ws.Range("E2").Value = "  -0.15%  "

$ws.Range("D3").Value = "1.826.35"
$ws.Range("E3").Value = "  -0.61%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9979"
$ws.Range("E4").Value = "  -0.27%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.57"
$ws.Range("E5").Value = "  -0.32%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6216"
$ws.Range("E6").Value = "  -1.16%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9998"
$ws.Range("E7").Value = "  -0.15%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07415"
$ws.Range("E8").Value = "  -2.05%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2915"
$ws.Range("E9").Value = "  -0.80%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.15"
$ws.Range("E10").Value = "  +2.36%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07686"
$ws.Range("E11").Value = "  -0.86%  "

$ws.Range("D12").Value = "1.836.28"
$ws.Range("E12").Value = "  -0.32%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.994"
$ws.Range("E13").Value = "  +0.51%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6656"
$ws.Range("E14").Value = "  -0.16%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "82.34"
$ws.Range("E15").Value = "  -1.02%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000009367"
$ws.Range("E16").Value = "  -6.75%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.935"
$ws.Range("E17").Value = "  -2.47%  "

$ws.Range("D18").Value = "29.058.12"
$ws.Range("E18").Value = "  -0.16%  "

$ws.Range("D19").Value = "2.072.31"
$ws.Range("E19").Value = "  -0.73%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.56"
$ws.Range("E20").Value = "  +0.95%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "222.56"
$ws.Range("E21").Value = "  -1.85%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  -0.15%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.096"
$ws.Range("E23").Value = "  -1.77%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9993"
$ws.Range("E24").Value = "  -0.25%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "159.61"
$ws.Range("E25").Value = "  -0.23%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1387"
$ws.Range("E26").Value = "  -0.32%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.464"
$ws.Range("E27").Value = "  -0.53%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.83"
$ws.Range("E28").Value = "  -0.63%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.488"
$ws.Range("E29").Value = "  -0.76%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05738"
$ws.Range("E30").Value = "  +8.96%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.148"
$ws.Range("E31").Value = "  +1.06%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.117"
$ws.Range("E32").Value = "  +2.44%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.209"

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.825"
$ws.Range("E34").Value = "  -1.32%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7367"
$ws.Range("E35").Value = "  -0.21%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.134"
$ws.Range("E36").Value = "  -0.31%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.664"
$ws.Range("E37").Value = "  -0.64%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.759"
$ws.Range("E38").Value = "  -0.17%  "

$ws.Range("D39").Value = "1.218.39"
$ws.Range("E39").Value = "  -2.15%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01767"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.470"
$ws.Range("E41").Value = "  +1.42%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8893"
$ws.Range("E42").Value = "  -1.60%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9991"
$ws.Range("E43").Value = "  -0.28%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.89"
$ws.Range("E44").Value = "  -0.20%  "

$ws.Range("D45").Value = "1.981.16"
$ws.Range("E45").Value = "  -0.36%  "

$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000124"
$ws.Range("E46").Value = "  -3.42%  "

$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "65.78"
$ws.Range("E47").Value = "  +2.16%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5077"
$ws.Range("E48").Value = "  -0.90%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07557"
$ws.Range("E49").Value = "  +14.83%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4037"
$ws.Range("E50").Value = "  -0.21%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.978"
$ws.Range("E51").Value = "  +0.30%  "
